$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hockey")
$ws.Activate()
$excel.ActiveWindow.TopLeftCell = "A23"
